$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant (used with PasteSpecial to copy only cell
# formatting/style from one range to another)
$xlPasteFormats = -4122

# --- 1) "Ativação:" date field: 01/01/2016 -> 01/01/2023 -------------------
# Force the cells to Text first so Excel doesn't auto-convert the
# date-shaped string into a real date serial number, then restore the
# original (General) formatting by pasting the format from a sibling
# row that already carries the correct style.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)

# --- 2) New "Objectives:" English text (row 11, B/C) ------------------------
$ws.Range("B11").Value = "Provide the student with the basic knowledge of electronic materials aiming their application in devices."
$ws.Range("C11").Value = "Provide the student with the basic knowledge of electronic materials aiming their application in devices."

$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)

# --- 3) New "Short syllabus:" English text (row 14, B/C) --------------------
$ws.Range("B14").Value = "Materials for electronics. Electronics and Solid State Physics. Semiconductor materials and devices. Optoelectronic materials and devices. Dielectric and piezoelectric materials and devices."
$ws.Range("C14").Value = "Materials for electronics. Electronics and Solid State Physics. Semiconductor materials and devices. Optoelectronic materials and devices. Dielectric and piezoelectric materials and devices."

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)

# --- 4) New "Syllabus:" English text (row 16, B/C) ---------------------------
$ws.Range("B16").Value = "Materials for electronic applications: metals, ceramics, glasses and polymers. Single crystals and thin films.Waves and particles in matter. Electrons in atoms and crystals. Energy band structures. Electronic and spectroscopic properties of materials.Conducting, semiconducting and insulating materials. Electronic properties in semiconductors. Electric transport. Semiconductor devices. pn junction Metal-semiconductor and semiconductor-insulator contact. Semiconductor devices: diodes and bipolar and FET transistors.Optoelectronic materials and devices. LED, semiconductor laser, photodetectors and photovoltaic cells.Types and properties of dielectric materials. Ferroelectric and piezoelectric materials. Devices based on dielectric and piezoelectric materials. Applications."
$ws.Range("C16").Value = "Materials for electronic applications: metals, ceramics, glasses and polymers. Single crystals and thin films.Waves and particles in matter. Electrons in atoms and crystals. Energy band structures. Electronic and spectroscopic properties of materials.Conducting, semiconducting and insulating materials. Electronic properties in semiconductors. Electric transport. Semiconductor devices. pn junction Metal-semiconductor and semiconductor-insulator contact. Semiconductor devices: diodes and bipolar and FET transistors.Optoelectronic materials and devices. LED, semiconductor laser, photodetectors and photovoltaic cells.Types and properties of dielectric materials. Ferroelectric and piezoelectric materials. Devices based on dielectric and piezoelectric materials. Applications."

$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
